# Updates cryptos list values (Price / Volume(1h) columns, plus two row swaps
# in B/C) to match the refreshed coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each written cell to stay plain text (NumberFormat "@") so values like
# "0.518" or "29.656.85" are not reinterpreted as numbers, then restore the
# "General" number format so formatting matches the original cells.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
}

Set-TextValue "D2" "29.656.85"
Set-TextValue "E2" "  +3.56%  "
Set-TextValue "D3" "1.604.86"
Set-TextValue "E3" "  +2.59%  "
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "212.25"
Set-TextValue "E5" "  +0.92%  "
Set-TextValue "D6" "0.518"
Set-TextValue "E6" "  +0.69%  "
Set-TextValue "E7" "  -0.07%  "
Set-TextValue "D8" "27.21"
Set-TextValue "E8" "  +9.50%  "
Set-TextValue "D9" "43.54"
Set-TextValue "E9" "  -1.33%  "
Set-TextValue "E10" "  +2.16%  "
Set-TextValue "E11" "  +2.13%  "
Set-TextValue "E12" "  +1.02%  "
Set-TextValue "D13" "1.835.58"
Set-TextValue "E13" "  +2.66%  "
Set-TextValue "D14" "1.602.71"
Set-TextValue "E14" "  +2.58%  "
Set-TextValue "D15" "29.665.65"
Set-TextValue "E15" "  +3.46%  "
Set-TextValue "E16" "  +3.82%  "
Set-TextValue "E17" "  +2.38%  "
Set-TextValue "D18" "63.34"
Set-TextValue "E18" "  +2.79%  "
Set-TextValue "D19" "241.00"
Set-TextValue "E19" "  +5.78%  "
Set-TextValue "E20" "  +3.83%  "
Set-TextValue "D21" "0.0₃0692"
Set-TextValue "E21" "  +1.47%  "
Set-TextValue "D22" "0.998"
Set-TextValue "E22" "  -0.11%  "
Set-TextValue "E23" "  +1.51%  "
Set-TextValue "D24" "9.22"
Set-TextValue "E24" "  +1.64%  "
Set-TextValue "D25" "2.08"
Set-TextValue "E25" "  +0.70%  "
Set-TextValue "D26" "155.19"
Set-TextValue "E26" "  +2.07%  "
Set-TextValue "D27" "15.35"
Set-TextValue "E27" "  +3.84%  "
Set-TextValue "E28" "  +0.87%  "
Set-TextValue "D29" "6.41"
Set-TextValue "E29" "  +2.55%  "
Set-TextValue "E30" "  -0.02%  "
Set-TextValue "E31" "  +3.72%  "
Set-TextValue "E32" "  +0.98%  "
Set-TextValue "E33" "  +1.10%  "
Set-TextValue "B34" "InternetComputer(DFINITY)"
Set-TextValue "C34" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "3.13"
Set-TextValue "E34" "  +4.08%  "
Set-TextValue "B35" "Maker"
Set-TextValue "C35" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D35" "1.430.93"
Set-TextValue "E35" "  +2.04%  "
Set-TextValue "E36" "  +0.21%  "
Set-TextValue "E37" "  +4.64%  "
Set-TextValue "D38" "2.82"
Set-TextValue "E38" "  +4.09%  "
Set-TextValue "E39" "  +0.18%  "
Set-TextValue "E41" "  +4.00%  "
Set-TextValue "E42" "  +1.66%  "
Set-TextValue "D43" "54.21"
Set-TextValue "E43" "  +27.45%  "
Set-TextValue "D44" "0.0487"
Set-TextValue "E44" "  +5.45%  "
Set-TextValue "D45" "0.800"
Set-TextValue "E45" "  +4.24%  "
Set-TextValue "E46" "  -0.12%  "
Set-TextValue "D47" "65.81"
Set-TextValue "E47" "  +2.92%  "
Set-TextValue "B48" "WEMIXToken"
Set-TextValue "C48" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D48" "0.946"
Set-TextValue "E48" "  +12.28%  "
Set-TextValue "B49" "FraxShare"
Set-TextValue "C49" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D49" "5.29"
Set-TextValue "E49" "  +1.35%  "
Set-TextValue "D50" "1.746.09"
Set-TextValue "E50" "  +2.86%  "
Set-TextValue "E51" "  +2.22%  "
